$wb = $excel.ActiveWorkbook

# Rename existing sheets
$wb.Worksheets.Item("SheetNew1").Name = "SheetNewNew1"
$wb.Worksheets.Item("SheetNew2").Name = "SheetNewNew2"

# Add new sheet "Sheet3" at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sheet3.Name = "Sheet3"

# Add new sheet "Sheet4" at the end (after Sheet3), this becomes active
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sheet4.Name = "Sheet4"

$sheet4.Activate()
